# Manual: New Section Task Switches. The general story, which does not yet
# explain the specifics of RTuinOS
#
# Adds a new row (52) to the "effort" sheet with:
#   A52 = 2012-12-03 (date, same number format as the rest of column A)
#   B52 = 2.5 (effort hours)
#   D52 = "Manual: New section "Task switches" started"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 52

$prevDateCell = $ws.Cells.Item($newRow - 1, 1)
$dateCell = $ws.Cells.Item($newRow, 1)
$prevDateCell.Copy()
$dateCell.PasteSpecial(-4122)  # xlPasteFormats
$dateCell.Value = Get-Date -Year 2012 -Month 12 -Day 3 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$ws.Cells.Item($newRow, 2).Value = 2.5

$ws.Cells.Item($newRow, 4).Value = 'Manual: New section "Task switches" started'

$ws.Range("E$newRow").Select()
